$wb = $excel.ActiveWorkbook

# --- Update the "Data" sheet values (capitalize Stock/Crypto) ---
$wsData = $wb.Worksheets.Item("Data")
$wsData.Range("A1").Value = "Stock"
$wsData.Range("A2").Value = "Crypto"

# --- Hide the "Data" sheet ---
$wsData.Visible = $False

# --- Update selection on Sheet1 ---
$wsSheet1 = $wb.Worksheets.Item("Sheet1")
$wsSheet1.Activate()
$wsSheet1.Range("B7").Select()
